# Normalization of the table in Dystrybucja database
# The IDDostawy key column (and its mirrored lookup column in the
# Ceny helper table) is converted from numeric identifiers (1-5) to
# text-style identifiers ("A1".."A5").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Ceny table (right-hand mini table), column N: IDCeny lookup ---
$ws.Range("N4").Value = "A1"
$ws.Range("N5").Value = "A2"
$ws.Range("N6").Value = "A3"
$ws.Range("N7").Value = "A4"
$ws.Range("N8").Value = "A5"

# --- Dostawy table, column B: IDDostawy primary key ---
$ws.Range("B12").Value = "A1"
$ws.Range("B13").Value = "A2"
$ws.Range("B14").Value = "A3"
$ws.Range("B15").Value = "A4"
$ws.Range("B16").Value = "A5"

# Update the active selection to match the author's final cursor position
$null = $ws.Range("B17").Select()
